{"js": "// Apply the LOQ4081.docx content revisions described by the commit diff.\n// Each entry is an exact, unique \"find -> replace\" pair for whole-text runs.\nconst replacements = [\n  [\n    \"Ativa\u00e7\u00e3o: 01/01/2018\",\n    \"Ativa\u00e7\u00e3o: 01/01/2022\",\n  ],\n  [\n    \"Propriedades f\u00edsicas (solubilidade, ponto de fus\u00e3o e ebuli\u00e7\u00e3o), acidez e basicidade dos compostos org\u00e2nicos. Classes de compostos org\u00e2nicos, principais rea\u00e7\u00f5es e introdu\u00e7\u00e3o \u00e0 mecanismos de rea\u00e7\u00f5es org\u00e2nicas (Substitui\u00e7\u00e3o, elimina\u00e7\u00e3o, adi\u00e7\u00e3o; via i\u00f4nica e radicalar). Import\u00e2ncia Industrial, impacto ambiental e degrada\u00e7\u00e3o de compostos org\u00e2nicos.\",\n    \"Propriedades f\u00edsicas, acidez e basicidade dos compostos org\u00e2nicos. Classes de compostos org\u00e2nicos, estereoquimica, principais rea\u00e7\u00f5es e introdu\u00e7\u00e3o \u00e0 mecanismos de rea\u00e7\u00f5es org\u00e2nicas (Substitui\u00e7\u00e3o, elimina\u00e7\u00e3o, adi\u00e7\u00e3o; via i\u00f4nica e radicalar). Import\u00e2ncia Industrial, impacto ambiental e degrada\u00e7\u00e3o de compostos org\u00e2nicos.\",\n  ],\n  [\n    \"Physical properties (solubility, melting point and boiling), acidity and basicity of organic compounds. Classes of organic compounds, main reactions and introduction to mechanisms of organic reactions (Substitution, elimination, addition, ionic and radical mechanism). Industrial Importance, environmental impact and degradation of organic compounds.\",\n    \"Physical properties (solubility, melting point and boiling), acidity and basicity of organic compounds. Classes of organic compounds, stereochemistry, main reactions and introduction to mechanisms of organic reactions (Substitution, elimination, addition, ionic and radical mechanism). Industrial Importance, environmental impact and degradation of organic compounds.\",\n  ],\n  [\n    \"Termologia e conceitos b\u00e1sicos de qu\u00edmica org\u00e2nica, estrutura das importantes classes dos compostos org\u00e2nicos, relacionando com produtos org\u00e2nicos antropog\u00eanicos. Discuss\u00e3o dos aspectos b\u00e1sicos dos processos de parti\u00e7\u00e3o. Como a estrutura qu\u00edmica determina a solubilidade do composto org\u00e2nico em \u00e1gua, presen\u00e7a de sais e co-solventes org\u00e2nicos. Principais classes de compostos org\u00e2nicos, principais rea\u00e7\u00f5es e introdu\u00e7\u00e3o \u00e0 mecanismos de rea\u00e7\u00f5es org\u00e2nicas (Substitui\u00e7\u00e3o, elimina\u00e7\u00e3o, adi\u00e7\u00e3o; via i\u00f4nica e radicalar). Como os absorventes naturais s\u00e3o importantes para transporte, distribui\u00e7\u00e3o e destino dos compostos org\u00e2nicos no meio ambiente. Processos de transforma\u00e7\u00e3o abi\u00f3tica e biol\u00f3gicas. Processos qu\u00edmicos, fotoqu\u00edmicos, biol\u00f3gicos e rea\u00e7\u00f5es de degrada\u00e7\u00e3o. Estudos de caso envolvendo os conceitos discutidos e sistemas ambientais, tais como lagos, rios e aqu\u00edferos.\",\n    \"Terminologia e conceitos b\u00e1sicos de qu\u00edmica org\u00e2nica, estrutura e caracter\u00edsticas das importantes classes dos compostos org\u00e2nicos, relacionando com produtos org\u00e2nicos antropog\u00eanicos. Acidez e basicidade em compostos org\u00e2nicos. Discuss\u00e3o dos aspectos b\u00e1sicos dos processos de parti\u00e7\u00e3o. Como a estrutura qu\u00edmica determina a solubilidade do composto org\u00e2nico em \u00e1gua. Estereoquimica. Principais rea\u00e7\u00f5es e introdu\u00e7\u00e3o \u00e0 mecanismos de rea\u00e7\u00f5es org\u00e2nicas (Substitui\u00e7\u00e3o, elimina\u00e7\u00e3o, adi\u00e7\u00e3o; via i\u00f4nica e radicalar). Estereoqu\u00edmica. Como os absorventes naturais s\u00e3o importantes para transporte, distribui\u00e7\u00e3o e destino dos compostos org\u00e2nicos no meio ambiente. Processos de transforma\u00e7\u00e3o abi\u00f3tica e biol\u00f3gicas. Processos qu\u00edmicos, fotoqu\u00edmicos, biol\u00f3gicos e rea\u00e7\u00f5es de degrada\u00e7\u00e3o. Estudos de caso envolvendo os conceitos discutidos e sistemas ambientais, tais como lagos, rios e aqu\u00edferos.\",\n  ],\n  [\n    \"Termology and basic concepts of organic chemistry, structure of important classes of organic compounds, relating to anthropogenic organic products. Discussion of the basics of partitioning processes. As the chemical structure determines the solubility of the organic compound in water, presence of salts and organic co-solvents. Main classes of organic compounds, main reactions and introduction to mechanisms of organic reactions (Substitution, elimination, addition, ionic and radical pathways). As natural absorbents are important for transport, distribution and fate of organic compounds in the environment. Abiotic and biological transformation processes. Chemical, photochemical, biological processes and degradation reactions. Case studies involving the concepts discussed and environmental systems, such as lakes, rivers and aquifers..\",\n    \"Terminology and basic concepts of organic chemistry, structure and characteristics of the important classes of organic compounds, relating to anthropogenic organic products. Acidity and basicity in organic compounds. Discussion of the basic aspects of partition processes. How the chemical structure determines the solubility of the organic compound in water. Stereochemistry. Main reactions and introduction to the organic reaction mechanisms (Substitution, elimination, addition; ionic and radical pathways). Stereochemistry. As natural absorbents are important for the transport, distribution and destination of organic compounds in the environment. Abiotic and biological transformation processes. Chemical, photochemical, biological processes and degradation reactions. Case studies involving the discussed concepts and environmental systems, such as lakes, rivers and aquifers.\",\n  ],\n  [\n    \"Bruice, Paula Yurkanis \u2013 Qu\u00edmica Org\u00e2nica \u2013 PEARSON Prentice Hall \u2013 S\u00e3o Paulo, 2006. Solomons, T.W.G.- Qu\u00edmica Org\u00e2nica. Volumes 1 e 2, Rio de Janeiro, Livros T\u00e9cnicos e Cient\u00edficos, 2012.Rene P. Schwarzenbach, Philip M. Gschwend, Dieter M. Imboden - Environmental Organic Chemistry \u2013 Wiley Interscience, 2016.\",\n    \"- Bruice, Paula Yurkanis \u2013 Qu\u00edmica Org\u00e2nica \u2013 PEARSON Prentice Hall \u2013 S\u00e3o Paulo, 2006. Solomons, T.W.G.- Qu\u00edmica Org\u00e2nica. Volumes 1 e 2, Rio de Janeiro, Livros T\u00e9cnicos e Cient\u00edficos, 2012.- Rene P. Schwarzenbach, Philip M. Gschwend, Dieter M. Imboden - Environmental Organic Chemistry \u2013 Wiley Interscience, 2016.- Madigan, M.T.; Martinko, J.M.; Bender, K.S.; Buckley, D.H.; Stahl, D.A. Microbiologia de Brock. Editora Artmed, 14a Edi\u00e7\u00e3o, 2016. - Nelson, D.; Cox, M. Princ\u00edpios de Bioqu\u00edmica de Lehninger. Artmed Editora. 6a  Edi\u00e7\u00e3o, 2014.- Pratt, C.; Cornely, K. Bioqu\u00edmica essencial. Guanabara Koogan. 1a  Edi\u00e7\u00e3o, 2006. - Wasserman, S.A.; Minorsky, P.V.; Jackson, R.; Reece, J.; Cain, M.; Urry, L. Biologia de Campbell. Artmed Editora. 8 a  Edi\u00e7\u00e3o. 2010.- Cooper, G.M. A C\u00e9lula \u2013 Uma Abordagem molecular. Artmed Editora Ltda. 3a  Edi\u00e7\u00e3o. 2007.- Raven, P.H.; Evert, S.E. Biologia vegetal. Editora Guanabara Koogan, 2007.- Maier, R. Environmental Microbiology. Academic Press. 2000. - Jordening, H.; Winter, J. Environmental Biotechnology. Concepts and Applications. Wiley-VCH. 2005. - Brock, T. D. ; Madigan, M.T.; Martinko, J.M.; Dunlap, P.; Clark, D. Biology of Microorganisms. Pearson Education.12a  Edi\u00e7\u00e3o. 2009.- Tortora, G.; Burdell, B.; Case, C. Microbiology. An Introduction. Pearson Benjamin Cummings. 10a  Edi\u00e7\u00e3o. 2010.\",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for \"${find.slice(0, 40)}...\" but found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(replace, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Apply the LOQ4081.docx content revisions described by the commit diff.\n# Each pair is an exact \"find -> replace\" whole-text match, applied with\n# Find.Execute (wdReplaceAll) against the full document Range.\n\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nfunction Replace-Exact($findText, $replaceText) {\n    $range = $d.Content\n    $ok = $range.Find.Execute(\n        $findText,      # FindText\n        $true,          # MatchCase\n        $false,         # MatchWholeWord\n        $false,         # MatchWildcards\n        $false,         # MatchSoundsLike\n        $false,         # MatchAllWordForms\n        $true,          # Forward\n        $wdFindContinue,# Wrap\n        $false,         # Format\n        $replaceText,   # ReplaceWith\n        $wdReplaceAll   # Replace\n    )\n    if (-not $ok) {\n        throw \"Find/replace failed for: $findText\"\n    }\n}\n\nReplace-Exact \"Ativa\u00e7\u00e3o: 01/01/2018\" \"Ativa\u00e7\u00e3o: 01/01/2022\"\n\nReplace-Exact \"Propriedades f\u00edsicas (solubilidade, ponto de fus\u00e3o e ebuli\u00e7\u00e3o), acidez e basicidade dos compostos org\u00e2nicos. Classes de compostos org\u00e2nicos, principais rea\u00e7\u00f5es e introdu\u00e7\u00e3o \u00e0 mecanismos de rea\u00e7\u00f5es org\u00e2nicas (Substitui\u00e7\u00e3o, elimina\u00e7\u00e3o, adi\u00e7\u00e3o; via i\u00f4nica e radicalar). Import\u00e2ncia Industrial, impacto ambiental e degrada\u00e7\u00e3o de compostos org\u00e2nicos.\" \"Propriedades f\u00edsicas, acidez e basicidade dos compostos org\u00e2nicos. Classes de compostos org\u00e2nicos, estereoquimica, principais rea\u00e7\u00f5es e introdu\u00e7\u00e3o \u00e0 mecanismos de rea\u00e7\u00f5es org\u00e2nicas (Substitui\u00e7\u00e3o, elimina\u00e7\u00e3o, adi\u00e7\u00e3o; via i\u00f4nica e radicalar). Import\u00e2ncia Industrial, impacto ambiental e degrada\u00e7\u00e3o de compostos org\u00e2nicos.\"\n\nReplace-Exact \"Physical properties (solubility, melting point and boiling), acidity and basicity of organic compounds. Classes of organic compounds, main reactions and introduction to mechanisms of organic reactions (Substitution, elimination, addition, ionic and radical mechanism). Industrial Importance, environmental impact and degradation of organic compounds.\" \"Physical properties (solubility, melting point and boiling), acidity and basicity of organic compounds. Classes of organic compounds, stereochemistry, main reactions and introduction to mechanisms of organic reactions (Substitution, elimination, addition, ionic and radical mechanism). Industrial Importance, environmental impact and degradation of organic compounds.\"\n\nReplace-Exact \"Termologia e conceitos b\u00e1sicos de qu\u00edmica org\u00e2nica, estrutura das importantes classes dos compostos org\u00e2nicos, relacionando com produtos org\u00e2nicos antropog\u00eanicos. Discuss\u00e3o dos aspectos b\u00e1sicos dos processos de parti\u00e7\u00e3o. Como a estrutura qu\u00edmica determina a solubilidade do composto org\u00e2nico em \u00e1gua, presen\u00e7a de sais e co-solventes org\u00e2nicos. Principais classes de compostos org\u00e2nicos, principais rea\u00e7\u00f5es e introdu\u00e7\u00e3o \u00e0 mecanismos de rea\u00e7\u00f5es org\u00e2nicas (Substitui\u00e7\u00e3o, elimina\u00e7\u00e3o, adi\u00e7\u00e3o; via i\u00f4nica e radicalar). Como os absorventes naturais s\u00e3o importantes para transporte, distribui\u00e7\u00e3o e destino dos compostos org\u00e2nicos no meio ambiente. Processos de transforma\u00e7\u00e3o abi\u00f3tica e biol\u00f3gicas. Processos qu\u00edmicos, fotoqu\u00edmicos, biol\u00f3gicos e rea\u00e7\u00f5es de degrada\u00e7\u00e3o. Estudos de caso envolvendo os conceitos discutidos e sistemas ambientais, tais como lagos, rios e aqu\u00edferos.\" \"Terminologia e conceitos b\u00e1sicos de qu\u00edmica org\u00e2nica, estrutura e caracter\u00edsticas das importantes classes dos compostos org\u00e2nicos, relacionando com produtos org\u00e2nicos antropog\u00eanicos. Acidez e basicidade em compostos org\u00e2nicos. Discuss\u00e3o dos aspectos b\u00e1sicos dos processos de parti\u00e7\u00e3o. Como a estrutura qu\u00edmica determina a solubilidade do composto org\u00e2nico em \u00e1gua. Estereoquimica. Principais rea\u00e7\u00f5es e introdu\u00e7\u00e3o \u00e0 mecanismos de rea\u00e7\u00f5es org\u00e2nicas (Substitui\u00e7\u00e3o, elimina\u00e7\u00e3o, adi\u00e7\u00e3o; via i\u00f4nica e radicalar). Estereoqu\u00edmica. Como os absorventes naturais s\u00e3o importantes para transporte, distribui\u00e7\u00e3o e destino dos compostos org\u00e2nicos no meio ambiente. Processos de transforma\u00e7\u00e3o abi\u00f3tica e biol\u00f3gicas. Processos qu\u00edmicos, fotoqu\u00edmicos, biol\u00f3gicos e rea\u00e7\u00f5es de degrada\u00e7\u00e3o. Estudos de caso envolvendo os conceitos discutidos e sistemas ambientais, tais como lagos, rios e aqu\u00edferos.\"\n\nReplace-Exact \"Termology and basic concepts of organic chemistry, structure of important classes of organic compounds, relating to anthropogenic organic products. Discussion of the basics of partitioning processes. As the chemical structure determines the solubility of the organic compound in water, presence of salts and organic co-solvents. Main classes of organic compounds, main reactions and introduction to mechanisms of organic reactions (Substitution, elimination, addition, ionic and radical pathways). As natural absorbents are important for transport, distribution and fate of organic compounds in the environment. Abiotic and biological transformation processes. Chemical, photochemical, biological processes and degradation reactions. Case studies involving the concepts discussed and environmental systems, such as lakes, rivers and aquifers..\" \"Terminology and basic concepts of organic chemistry, structure and characteristics of the important classes of organic compounds, relating to anthropogenic organic products. Acidity and basicity in organic compounds. Discussion of the basic aspects of partition processes. How the chemical structure determines the solubility of the organic compound in water. Stereochemistry. Main reactions and introduction to the organic reaction mechanisms (Substitution, elimination, addition; ionic and radical pathways). Stereochemistry. As natural absorbents are important for the transport, distribution and destination of organic compounds in the environment. Abiotic and biological transformation processes. Chemical, photochemical, biological processes and degradation reactions. Case studies involving the discussed concepts and environmental systems, such as lakes, rivers and aquifers.\"\n\nReplace-Exact \"Bruice, Paula Yurkanis \u2013 Qu\u00edmica Org\u00e2nica \u2013 PEARSON Prentice Hall \u2013 S\u00e3o Paulo, 2006. Solomons, T.W.G.- Qu\u00edmica Org\u00e2nica. Volumes 1 e 2, Rio de Janeiro, Livros T\u00e9cnicos e Cient\u00edficos, 2012.Rene P. Schwarzenbach, Philip M. Gschwend, Dieter M. Imboden - Environmental Organic Chemistry \u2013 Wiley Interscience, 2016.\" \"- Bruice, Paula Yurkanis \u2013 Qu\u00edmica Org\u00e2nica \u2013 PEARSON Prentice Hall \u2013 S\u00e3o Paulo, 2006. Solomons, T.W.G.- Qu\u00edmica Org\u00e2nica. Volumes 1 e 2, Rio de Janeiro, Livros T\u00e9cnicos e Cient\u00edficos, 2012.- Rene P. Schwarzenbach, Philip M. Gschwend, Dieter M. Imboden - Environmental Organic Chemistry \u2013 Wiley Interscience, 2016.- Madigan, M.T.; Martinko, J.M.; Bender, K.S.; Buckley, D.H.; Stahl, D.A. Microbiologia de Brock. Editora Artmed, 14a Edi\u00e7\u00e3o, 2016. - Nelson, D.; Cox, M. Princ\u00edpios de Bioqu\u00edmica de Lehninger. Artmed Editora. 6a  Edi\u00e7\u00e3o, 2014.- Pratt, C.; Cornely, K. Bioqu\u00edmica essencial. Guanabara Koogan. 1a  Edi\u00e7\u00e3o, 2006. - Wasserman, S.A.; Minorsky, P.V.; Jackson, R.; Reece, J.; Cain, M.; Urry, L. Biologia de Campbell. Artmed Editora. 8 a  Edi\u00e7\u00e3o. 2010.- Cooper, G.M. A C\u00e9lula \u2013 Uma Abordagem molecular. Artmed Editora Ltda. 3a  Edi\u00e7\u00e3o. 2007.- Raven, P.H.; Evert, S.E. Biologia vegetal. Editora Guanabara Koogan, 2007.- Maier, R. Environmental Microbiology. Academic Press. 2000. - Jordening, H.; Winter, J. Environmental Biotechnology. Concepts and Applications. Wiley-VCH. 2005. - Brock, T. D. ; Madigan, M.T.; Martinko, J.M.; Dunlap, P.; Clark, D. Biology of Microorganisms. Pearson Education.12a  Edi\u00e7\u00e3o. 2009.- Tortora, G.; Burdell, B.; Case, C. Microbiology. An Introduction. Pearson Benjamin Cummings. 10a  Edi\u00e7\u00e3o. 2010.\"\n"}
